$wb = $excel.ActiveWorkbook

# --- Create the three new sheets in order, after Sheet1 ---
$wsReport = $wb.Worksheets.Add()
$wsReport.Name = "DoAiReport Prompts"
$wb.Worksheets("DoAiReport Prompts").Move($null, $wb.Worksheets("Sheet1"))

$wsPy = $wb.Worksheets.Add()
$wsPy.Name = "DoAiReportPy Prompts"
$wb.Worksheets("DoAiReportPy Prompts").Move($null, $wb.Worksheets("DoAiReport Prompts"))

$wsHtml = $wb.Worksheets.Add()
$wsHtml.Name = "DoAiReportHtml Prompts"
$wb.Worksheets("DoAiReportHtml Prompts").Move($null, $wb.Worksheets("DoAiReportPy Prompts"))

# --- Sheet2: DoAiReport Prompts ---
$ws2 = $wb.Worksheets("DoAiReport Prompts")
$ws2.Range("A1").Value = "DoAiReport.query"
$ws2.Range("A2").Value = "Build a table describing failed test cases."
$ws2.Range("A3").Value = "Create a report describing test case failures. Format as table. Group failures of the same test case and reduce the number of lines in the table for same errors."
$ws2.Range("A4").Value = "Describe any patterns and regularities you can see in this data."
$ws2.Range("A5").Value = "Find groups of test cases that usually fail together."
$ws2.Range("A6").Value = "What are the top 10 longest test runs?"
$ws2.Range("A7").Value = "Do you observe any anomalies in this test run data?"
$ws2.Range("A8").Value = "Build a list of failed test cases and their fail reasons."
$ws2.Range("A9").Value = "Create a report for this sequence of test runs. Include failures only into a summary table. Also provide analysis for all runs and list unique TestSet names."
$ws2.Range("A1").Font.Bold = $true
$ws2.Columns("A").ColumnWidth = 150.96354166666666

# --- Sheet3: DoAiReportPy Prompts ---
$ws3 = $wb.Worksheets("DoAiReportPy Prompts")
$ws3.Range("A1").Value = "DoAiReportPy.query"
$ws3.Range("A2").Value = "Draw status for each test case over time. Each test case - is a horizontal sequence of dots."
$ws3.Range("A3").Value = "Discover which test cases have at least one non-pass status (build an array). Go again through all data and capture available runs of these test cases. Draw all statuses (including Pass, etc.) for each such test case over time. Each test case - is a horizontal sequence of dots. "
$ws3.Range("A4").Value = "Draw fail status for each test case over time. Each test case - is a horizontal sequence of dots. Do not draw Pass points."
$ws3.Range("A1").Font.Bold = $true
$ws3.Columns("A").ColumnWidth = 144.66276041666666

# --- Sheet4: DoAiReportHtml Prompts ---
$ws4 = $wb.Worksheets("DoAiReportHtml Prompts")
$ws4.Range("A1").Value = "DoAiReportHtml.query"
$ws4.Range("A2").Value = "Discover which test cases have at least one non-pass status (build an array). Go again through all data and capture available runs of these test cases. Draw all statuses (including Pass, etc.) for each such test case over time. Each test case - is a horizontal sequence of dots."
$ws4.Range("A1").Font.Bold = $true
$ws4.Columns("A").ColumnWidth = 114.46354166666667

# --- Select Sheet4 (DoAiReportHtml Prompts) as the active/visible sheet (tabSelected) ---
$ws4.Activate()
